$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 4289.854091434014
$ws.Range("C3").Value = 4272.842298340157
$ws.Range("C4").Value = 4211.445205727659
$ws.Range("C5").Value = 4211.445205727659
$ws.Range("C6").Value = 4211.445205727659
$ws.Range("C7").Value = 4211.445205727659
$ws.Range("C8").Value = 4090.793752666073
$ws.Range("C9").Value = 4090.793752666073
$ws.Range("C10").Value = 4090.793752666073
$ws.Range("C11").Value = 4077.085650267751
$ws.Range("C12").Value = 4077.085650267751
